# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
# All values in this workbook are stored as plain text (dates like
# "2026-01-28", times like "18:13:53"/"18:00", percentages like "88.2%"
# would otherwise be auto-converted by Excel into dates/numbers), so for
# each new range we force Text number format before writing the values and
# then clear the formatting again afterwards so the appended rows end up
# with the same "no explicit style" look as the rest of the sheet.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        $SheetName,
        $StartRow,
        $Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $endRow = $StartRow + $Rows.Count - 1
    $rng = $ws.Range("A$StartRow`:F$endRow")

    # Force text formatting first so date/time/percentage-looking strings
    # are not auto-converted into numbers/dates.
    $rng.NumberFormat = "@"

    $r = $StartRow
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r = $r + 1
    }

    # Drop the explicit "Text" style again so the new cells match the
    # unstyled look of the existing rows.
    $rng.ClearFormats()
}

# ---- PIR sheet: rows 150-162 ----
$pirRows = @(
    ,@("2026-01-28", "18:13:53", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:13:55", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:13:58", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:01", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:06", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:11", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:17", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:21", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:26", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:31", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:37", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:41", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-01-28", "18:14:46", "18:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows "PIR" 150 $pirRows

# ---- Humidity sheet: rows 144-156 ----
$humidityRows = @(
    ,@("2026-01-28", "18:13:51", "18:00", "Bathroom", "88.2%", "Active")
    ,@("2026-01-28", "18:13:53", "18:00", "Bathroom", "87.3%", "Active")
    ,@("2026-01-28", "18:13:56", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:13:59", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:07", "18:00", "Bathroom", "87.3%", "Active")
    ,@("2026-01-28", "18:14:12", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:16", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:24", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:28", "18:00", "Bathroom", "87.4%", "Active")
    ,@("2026-01-28", "18:14:32", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:36", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:44", "18:00", "Bathroom", "88.3%", "Active")
    ,@("2026-01-28", "18:14:48", "18:00", "Bathroom", "87.3%", "Active")
)
Add-LogRows "Humidity" 144 $humidityRows

# ---- Temperature sheet: rows 144-156 ----
$temperatureRows = @(
    ,@("2026-01-28", "18:13:52", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:13:54", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:13:57", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:00", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:08", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:13", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:16", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:25", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:28", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:33", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:37", "18:00", "Bathroom", "23.0C", "Active")
    ,@("2026-01-28", "18:14:45", "18:00", "Bathroom", "22.9C", "Active")
    ,@("2026-01-28", "18:14:49", "18:00", "Bathroom", "22.9C", "Active")
)
Add-LogRows "Temperature" 144 $temperatureRows
